$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.78559999999999
$ws.Range("D3").Value = -5.867699999999998
$ws.Range("E5").Value = 12.7039
$ws.Range("D14").Value = -7.807799999999998
$ws.Range("D21").Value = -7.669799999999997
$ws.Range("D23").Value = -7.182999999999992
$ws.Range("D25").Value = -8.372700000000002
